$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.212.17"
$ws.Range("E2").Value = "  -2.54%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.804.73"
$ws.Range("E3").Value = "  -0.55%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.69%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.33"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.63%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5307"
$ws.Range("E7").Value = "  -2.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3784"
$ws.Range("E8").Value = "  -1.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07483"
$ws.Range("E9").Value = "  -1.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.99"
$ws.Range("E10").Value = "  -1.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.099"
$ws.Range("E11").Value = "  -3.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  +0.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.213"
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.59"
$ws.Range("E14").Value = "  -3.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.340"
$ws.Range("E15").Value = "  -1.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.802.83"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.50"
$ws.Range("E17").Value = "  -2.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001069"
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06507"
$ws.Range("E19").Value = "  +1.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.24"
$ws.Range("E21").Value = "  -1.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.923"
$ws.Range("E22").Value = "  -1.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.225.02"
$ws.Range("E23").Value = "  -2.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.18"
$ws.Range("E24").Value = "  -2.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.089"
$ws.Range("E25").Value = "  -1.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.52"
$ws.Range("E26").Value = "  -4.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.49"
$ws.Range("E27").Value = "  -1.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.010.25"
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.331"
$ws.Range("E29").Value = "  -4.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "122.06"
$ws.Range("E30").Value = "  -1.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.125"
$ws.Range("E31").Value = "  -2.96%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1102"
$ws.Range("E32").Value = "  +7.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.590"
$ws.Range("E33").Value = "  -4.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.615"
$ws.Range("E34").Value = "  -1.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07222"
$ws.Range("E35").Value = "  +9.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2228"
$ws.Range("E36").Value = "  -4.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02299"
$ws.Range("E37").Value = "  -1.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.072"
$ws.Range("E38").Value = "  -1.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.530"
$ws.Range("E39").Value = "  -1.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6171"
$ws.Range("E40").Value = "  -3.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.13"
$ws.Range("E41").Value = "  -5.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.184"
$ws.Range("E42").Value = "  -4.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.433"
$ws.Range("E43").Value = "  +2.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.37"
$ws.Range("E45").Value = "  -2.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.681"
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5770"
$ws.Range("E47").Value = "  -4.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.68"
$ws.Range("E48").Value = "  -0.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.195"
$ws.Range("E49").Value = "  +2.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.927"
$ws.Range("E50").Value = "  -4.56%  "
$ws.Range("E51").Value = "  -2.41%  "
